# Update KHL probabilities tour sheet with the 2025-12-18 data refresh:
# overwrite rows 2-3 (replaced games) and append new rows 4-7 (newly scraped games).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Row 2: Автомобилист vs Динамо Мн
$ws.Range("A2").Value = 1369
$ws.Range("B2").Value = "2025-12-18T17:00:00"
$ws.Range("C2").Value = "Автомобилист"
$ws.Range("D2").Value = "Динамо Мн"
$ws.Range("E2").Value = 897879
$ws.Range("F2").Value = "https://text.khl.ru/text/897879.html"
$ws.Range("G2").Value = 4.725
$ws.Range("H2").Value = 3.748293
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 1.48
$ws.Range("K2").Value = 3.1025
$ws.Range("L2").Value = 3.624147
$ws.Range("M2").Value = 8.473293
$ws.Range("N2").Value = 31.836521
$ws.Range("O2").Value = 36.171304
$ws.Range("P2").Value = 68.007825
$ws.Range("Q2").Value = 0.16
$ws.Range("R2").Value = 0.099772
$ws.Range("S2").Value = 0.344335
$ws.Range("T2").Value = 0.154061
$ws.Range("U2").Value = 0.49988
$ws.Range("V2").Value = 0.097176
$ws.Range("W2").Value = 0.901099
$ws.Range("X2").Value = 0.19942
$ws.Range("Y2").Value = 0.798856
$ws.Range("Z2").Value = 0.336972
$ws.Range("AA2").Value = 0.661304
$ws.Range("AB2").Value = 0.491182
$ws.Range("AC2").Value = 0.507093
$ws.Range("AD2").Value = 0.63937
$ws.Range("AE2").Value = 0.358905
$ws.Range("AF2").Value = 0.815647
$ws.Range("AG2").Value = 0.184353
$ws.Range("AH2").Value = 0.599378
$ws.Range("AI2").Value = 0.400622
$ws.Range("AJ2").Value = 0.876665
$ws.Range("AK2").Value = 0.123335
$ws.Range("AL2").Value = 0.701506
$ws.Range("AM2").Value = 0.298494
$ws.Range("AN2").Value = 0.651934
$ws.Range("AO2").Value = 0.785379

# Row 3: Нефтехимик vs Авангард
$ws.Range("A3").Value = 1369
$ws.Range("B3").Value = "2025-12-18T19:00:00"
$ws.Range("C3").Value = "Нефтехимик"
$ws.Range("D3").Value = "Авангард"
$ws.Range("E3").Value = 897876
$ws.Range("F3").Value = "https://text.khl.ru/text/897876.html"
$ws.Range("G3").Value = 3.591304
$ws.Range("H3").Value = 5.075
$ws.Range("I3").Value = 4.504348
$ws.Range("J3").Value = 3.9375
$ws.Range("K3").Value = 3.764402
$ws.Range("L3").Value = 4.789674
$ws.Range("M3").Value = 8.666304
$ws.Range("N3").Value = 34.740952
$ws.Range("O3").Value = 39.952501
$ws.Range("P3").Value = 74.693453
$ws.Range("Q3").Value = 0.16
$ws.Range("R3").Value = 0.16
$ws.Range("S3").Value = 0.297553
$ws.Range("T3").Value = 0.13074
$ws.Range("U3").Value = 0.559664
$ws.Range("V3").Value = 0.029002
$ws.Range("W3").Value = 0.958955
$ws.Range("X3").Value = 0.072005
$ws.Range("Y3").Value = 0.915953
$ws.Range("Z3").Value = 0.145574
$ws.Range("AA3").Value = 0.842384
$ws.Range("AB3").Value = 0.25046
$ws.Range("AC3").Value = 0.737498
$ws.Range("AD3").Value = 0.378632
$ws.Range("AE3").Value = 0.609326
$ws.Range("AF3").Value = 0.889554
$ws.Range("AG3").Value = 0.110446
$ws.Range("AH3").Value = 0.725305
$ws.Range("AI3").Value = 0.274695
$ws.Range("AJ3").Value = 0.951858
$ws.Range("AK3").Value = 0.048142
$ws.Range("AL3").Value = 0.856479
$ws.Range("AM3").Value = 0.143521
$ws.Range("AN3").Value = 0.566774
$ws.Range("AO3").Value = 0.799243

# Row 4: Локомотив vs ЦСКА
$ws.Range("A4").Value = 1369
$ws.Range("B4").Value = "2025-12-18T19:00:00"
$ws.Range("C4").Value = "Локомотив"
$ws.Range("D4").Value = "ЦСКА"
$ws.Range("E4").Value = 897877
$ws.Range("F4").Value = "https://text.khl.ru/text/897877.html"
$ws.Range("G4").Value = 2.852452
$ws.Range("H4").Value = 2.846302
$ws.Range("I4").Value = 2.373399
$ws.Range("J4").Value = 1.610526
$ws.Range("K4").Value = 2.231489
$ws.Range("L4").Value = 2.60985
$ws.Range("M4").Value = 5.698754
$ws.Range("N4").Value = 27.320471
$ws.Range("O4").Value = 27.63696
$ws.Range("P4").Value = 54.957431
$ws.Range("Q4").Value = -0.051904
$ws.Range("R4").Value = 0.080709
$ws.Range("S4").Value = 0.34135
$ws.Range("T4").Value = 0.184285
$ws.Range("U4").Value = 0.474253
$ws.Range("V4").Value = 0.288008
$ws.Range("W4").Value = 0.71188
$ws.Range("X4").Value = 0.46876
$ws.Range("Y4").Value = 0.531128
$ws.Range("Z4").Value = 0.643777
$ws.Range("AA4").Value = 0.356111
$ws.Range("AB4").Value = 0.784996
$ws.Range("AC4").Value = 0.214892
$ws.Range("AD4").Value = 0.882666
$ws.Range("AE4").Value = 0.117222
$ws.Range("AF4").Value = 0.65304
$ws.Range("AG4").Value = 0.34696
$ws.Range("AH4").Value = 0.385717
$ws.Range("AI4").Value = 0.614283
$ws.Range("AJ4").Value = 0.734512
$ws.Range("AK4").Value = 0.265488
$ws.Range("AL4").Value = 0.484041
$ws.Range("AM4").Value = 0.515959
$ws.Range("AN4").Value = 0.702849
$ws.Range("AO4").Value = 0.810061

# Row 5: Северсталь vs Салават Юлаев
$ws.Range("A5").Value = 1369
$ws.Range("B5").Value = "2025-12-18T19:00:00"
$ws.Range("C5").Value = "Северсталь"
$ws.Range("D5").Value = "Салават Юлаев"
$ws.Range("E5").Value = 897878
$ws.Range("F5").Value = "https://text.khl.ru/text/897878.html"
$ws.Range("G5").Value = 1.666667
$ws.Range("H5").Value = 1.285714
$ws.Range("I5").Value = 1.333333
$ws.Range("J5").Value = 3.866667
$ws.Range("K5").Value = 2.766667
$ws.Range("L5").Value = 1.309524
$ws.Range("M5").Value = 2.952381
$ws.Range("N5").Value = 24.837331
$ws.Range("O5").Value = 22.434051
$ws.Range("P5").Value = 47.271382
$ws.Range("Q5").Value = -0.16
$ws.Range("R5").Value = -0.16
$ws.Range("S5").Value = 0.679401
$ws.Range("T5").Value = 0.162469
$ws.Range("U5").Value = 0.157982
$ws.Range("V5").Value = 0.418729
$ws.Range("W5").Value = 0.581123
$ws.Range("X5").Value = 0.613955
$ws.Range("Y5").Value = 0.385897
$ws.Range("Z5").Value = 0.773112
$ws.Range("AA5").Value = 0.22674
$ws.Range("AB5").Value = 0.881237
$ws.Range("AC5").Value = 0.118615
$ws.Range("AD5").Value = 0.9442
$ws.Range("AE5").Value = 0.055652
$ws.Range("AF5").Value = 0.763185
$ws.Range("AG5").Value = 0.236815
$ws.Range("AH5").Value = 0.522563
$ws.Range("AI5").Value = 0.477437
$ws.Range("AJ5").Value = 0.376547
$ws.Range("AK5").Value = 0.623453
$ws.Range("AL5").Value = 0.145086
$ws.Range("AM5").Value = 0.854914
$ws.Range("AN5").Value = 0.937513
$ws.Range("AO5").Value = 0.522518

# Row 6: Динамо М vs Барыс
$ws.Range("A6").Value = 1369
$ws.Range("B6").Value = "2025-12-18T19:30:00"
$ws.Range("C6").Value = "Динамо М"
$ws.Range("D6").Value = "Барыс"
$ws.Range("E6").Value = 897875
$ws.Range("F6").Value = "https://text.khl.ru/text/897875.html"
$ws.Range("G6").Value = 2.307876
$ws.Range("H6").Value = 1.32
$ws.Range("I6").Value = 1.376471
$ws.Range("J6").Value = 4.48
$ws.Range("K6").Value = 3.393938
$ws.Range("L6").Value = 1.348235
$ws.Range("M6").Value = 3.627876
$ws.Range("N6").Value = 28.935124
$ws.Range("O6").Value = 24.347309
$ws.Range("P6").Value = 53.282433
$ws.Range("Q6").Value = -0.051254
$ws.Range("R6").Value = -0.16
$ws.Range("S6").Value = 0.761412
$ws.Range("T6").Value = 0.125491
$ws.Range("U6").Value = 0.112298
$ws.Range("V6").Value = 0.303097
$ws.Range("W6").Value = 0.696104
$ws.Range("X6").Value = 0.486835
$ws.Range("Y6").Value = 0.512366
$ws.Range("Z6").Value = 0.661098
$ws.Range("AA6").Value = 0.338103
$ws.Range("AB6").Value = 0.798829
$ws.Range("AC6").Value = 0.200372
$ws.Range("AD6").Value = 0.892136
$ws.Range("AE6").Value = 0.107065
$ws.Range("AF6").Value = 0.852468
$ws.Range("AG6").Value = 0.147532
$ws.Range("AH6").Value = 0.659089
$ws.Range("AI6").Value = 0.340911
$ws.Range("AJ6").Value = 0.390168
$ws.Range("AK6").Value = 0.609832
$ws.Range("AL6").Value = 0.154136
$ws.Range("AM6").Value = 0.845864
$ws.Range("AN6").Value = 0.955986
$ws.Range("AO6").Value = 0.411692

# Row 7: СКА vs Ак Барс
$ws.Range("A7").Value = 1369
$ws.Range("B7").Value = "2025-12-18T19:30:00"
$ws.Range("C7").Value = "СКА"
$ws.Range("D7").Value = "Ак Барс"
$ws.Range("E7").Value = 897893
$ws.Range("F7").Value = "https://text.khl.ru/text/897893.html"
$ws.Range("G7").Value = 3.952203
$ws.Range("H7").Value = 2.02599
$ws.Range("I7").Value = 1.895858
$ws.Range("J7").Value = 1.325136
$ws.Range("K7").Value = 2.63867
$ws.Range("L7").Value = 1.960924
$ws.Range("M7").Value = 5.978193
$ws.Range("N7").Value = 32.220977
$ws.Range("O7").Value = 30.349418
$ws.Range("P7").Value = 62.570395
$ws.Range("Q7").Value = 0.097685
$ws.Range("R7").Value = -0.078892
$ws.Range("S7").Value = 0.52897
$ws.Range("T7").Value = 0.183607
$ws.Range("U7").Value = 0.287318
$ws.Range("V7").Value = 0.325773
$ws.Range("W7").Value = 0.674122
$ws.Range("X7").Value = 0.51331
$ws.Range("Y7").Value = 0.486584
$ws.Range("Z7").Value = 0.68583
$ws.Range("AA7").Value = 0.314065
$ws.Range("AB7").Value = 0.818083
$ws.Range("AC7").Value = 0.181812
$ws.Range("AD7").Value = 0.904984
$ws.Range("AE7").Value = 0.09491
$ws.Range("AF7").Value = 0.739994
$ws.Range("AG7").Value = 0.260006
$ws.Range("AH7").Value = 0.491234
$ws.Range("AI7").Value = 0.508766
$ws.Range("AJ7").Value = 0.583314
$ws.Range("AK7").Value = 0.416686
$ws.Range("AL7").Value = 0.312748
$ws.Range("AM7").Value = 0.687252
$ws.Range("AN7").Value = 0.852143
$ws.Range("AO7").Value = 0.658729

